# Settings back to qfin22 presentation - recalibration
# Updates the autocorrelation-of-abs-residuals arrays (shared strings) and
# the recalibrated parameter values on each model sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# AR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = -0.02767448908363682
$ws.Range("B3").Value = 0.8047757127757948
$ws.Range("B4").Value = 43.11434456897841
$ws.Range("B5").Value = "[1.0, 0.19426876504229537, 0.16364120572599727, 0.17585724837010477, 0.13823781080373698, 0.36997010207184206, 0.14746345388066803, 0.15781724496043834, 0.16439551323188364, 0.1005243254489963, 0.10499708069771696, 0.10641108231727008, 0.10554637556299769, 0.12828189141690693, 0.1476635650893182, 0.12140421076836844, 0.1109614396113983, 0.08022625785982869, 0.09399171699438791, 0.08956071036248883]"

# ---------------------------------------------------------------------------
# SETAR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -2.437708426073025
$ws.Range("B5").Value = 0.7115725858861324
$ws.Range("B6").Value = 31.34385888964406
$ws.Range("B7").Value = 2.355282900630713
$ws.Range("B8").Value = 0.69230112966093
$ws.Range("B9").Value = 27.06559673119342
$ws.Range("B10").Value = "[1.0, 0.22034013690219126, 0.19387650748222593, 0.19030872883139097, 0.19432005671143426, 0.24496777617322266, 0.12460059693315556, 0.15444955974900518, 0.14858492367430082, 0.11790317505964992, 0.13412707583333314, 0.11555686209194008, 0.127491440693676, 0.13508331042851668, 0.15290803071864775, 0.10053001921307084, 0.10575958781562801, 0.10673933689254754, 0.10531566068994845, 0.08189275056326766]"

# ---------------------------------------------------------------------------
# GARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = -0.02312415407085303
$ws.Range("B3").Value = 0.8589133746316927
$ws.Range("B4").Value = 0.07307536315403911
$ws.Range("B5").Value = 0.9093282917998065
$ws.Range("B6").Value = "[1.0, -0.008942395604581113, -0.03115342952507032, -0.025864828039431704, -0.06276178704619029, 0.3062341198810408, -0.01517746074436223, 0.01065377428082957, 0.005570344630911677, -0.038735829245318, -0.04540306599321187, -0.05615598669290948, -0.04199200721831009, -0.01286598315480665, 0.0333257910731646, -0.006495119450522699, -0.01117146604527087, -0.0476798266196217, -0.052773188307566256, -0.01870700524106385]"

# ---------------------------------------------------------------------------
# TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.05076866730824595
$ws.Range("B3").Value = 0.8552331217361768
$ws.Range("B4").Value = 0.06536518288929777
$ws.Range("B5").Value = 0.9101970434274161
$ws.Range("B6").Value = "[1.0, -0.008519285975574479, -0.0309165188057894, -0.02551564631802768, -0.062114053872304925, 0.30603701410650397, -0.014617398254979184, 0.010968026493721444, 0.006025128634411305, -0.03815160139853932, -0.044780929655319754, -0.055931759287914565, -0.04177466653874382, -0.012862621113906991, 0.033412863055942224, -0.006447284614206047, -0.010940810679177822, -0.04736295102706475, -0.052328517382627816, -0.018856218139420764]"
$ws.Range("B7").Value = 0.01366213949475307

# ---------------------------------------------------------------------------
# AR_TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = -0.0433713319745662
$ws.Range("B3").Value = 0.8047827018790829
$ws.Range("B4").Value = 0.0812875597365546
$ws.Range("B5").Value = 0.9059948725649249
$ws.Range("B6").Value = "[1.0, 0.011593150842141655, -0.003633802070473796, 0.001247896946016283, -0.029875748917702012, 0.2340350392672941, -0.017155907574380166, 0.00489354943360922, -0.002407462662252031, -0.03718396362152853, -0.043037766435038496, -0.042974635944310176, -0.028715440626545168, -0.019670124000556406, 0.01929272757528464, -0.001555355299112527, -0.010684511205526196, -0.029936365976336633, -0.03939273853083118, -0.033495783887166855]"
$ws.Range("B7").Value = -0.01133801560033445
$ws.Range("B9").Value = 0.7934167970366459
